$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vize")
$ws.Range("A1").Value = 2
$ws.Range("F3").Value = 2
